$d = $word.ActiveDocument
$d.Content.Find.Execute("Operating systems on EC2 instances", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Operating systems on Amazon EC2 instances", 2)
